$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("D2").Value = 74
$ws.Range("I2").Value = 93
$ws.Range("J2").Value = 91
$ws.Range("D3").Value = 106
$ws.Range("G3").Value = 107
$ws.Range("H3").Value = 108
$ws.Range("J3").Value = 173
$ws.Range("C9").Value = 374
$ws.Range("D9").Value = 328
$ws.Range("E9").Value = 347
$ws.Range("F9").Value = 410
$ws.Range("G9").Value = 381
$ws.Range("H9").Value = 359
$ws.Range("J9").Value = 325
$ws.Range("B10").Value = 1028
$ws.Range("C10").Value = 1235
$ws.Range("D10").Value = 1405
$ws.Range("E10").Value = 1711
$ws.Range("F10").Value = 1734
$ws.Range("H10").Value = 449
$ws.Range("I10").Value = 694
$ws.Range("J10").Value = 566
$ws.Range("K10").Value = 559
$ws.Range("B11").Value = 1438
$ws.Range("C11").Value = 1741
$ws.Range("D11").Value = 1923
$ws.Range("E11").Value = 2231
$ws.Range("F11").Value = 2320
$ws.Range("G11").Value = 1360
$ws.Range("H11").Value = 1020
$ws.Range("I11").Value = 1388
$ws.Range("J11").Value = 1179
$ws.Range("K11").Value = 1279

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("B2").Value = 7
$ws.Range("B7").Value = 42
$ws.Range("F8").Value = 107
$ws.Range("H8").Value = 81
$ws.Range("C19").Value = 44
$ws.Range("G20").Value = 15
$ws.Range("D22").Value = 15
$ws.Range("J23").Value = 10
$ws.Range("C28").Value = 113
$ws.Range("F28").Value = 93
$ws.Range("G28").Value = 70
$ws.Range("I28").Value = 75
$ws.Range("J28").Value = 44
$ws.Range("F32").Value = 159
$ws.Range("J32").Value = 62
$ws.Range("D36").Value = 66
$ws.Range("J42").Value = 14
$ws.Range("B45").Value = 22
$ws.Range("B47").Value = 44
$ws.Range("E47").Value = 58
$ws.Range("H47").Value = 33
$ws.Range("F49").Value = 9
$ws.Range("F50").Value = 52
$ws.Range("I50").Value = 18
$ws.Range("D52").Value = 31
$ws.Range("B53").Value = 187
$ws.Range("C53").Value = 278
$ws.Range("D53").Value = 472
$ws.Range("E53").Value = 559
$ws.Range("F53").Value = 524
$ws.Range("G53").Value = 213
$ws.Range("J53").Value = 184
$ws.Range("K53").Value = 167
$ws.Range("C55").Value = 6
$ws.Range("D61").Value = 21
$ws.Range("F61").Value = 36
$ws.Range("D62").Value = 21
$ws.Range("C65").Value = 47
$ws.Range("J68").Value = 12
$ws.Range("I74").Value = 37
$ws.Range("C76").Value = 61
$ws.Range("K76").Value = 40
$ws.Range("B77").Value = 64
$ws.Range("F78").Value = 42
$ws.Range("H80").Value = 13
$ws.Range("F83").Value = 23
$ws.Range("C87").Value = 34
$ws.Range("E95").Value = 73
$ws.Range("B99").Value = 1438
$ws.Range("C99").Value = 1741
$ws.Range("D99").Value = 1923
$ws.Range("E99").Value = 2231
$ws.Range("F99").Value = 2320
$ws.Range("G99").Value = 1360
$ws.Range("H99").Value = 1020
$ws.Range("I99").Value = 1388
$ws.Range("J99").Value = 1179
$ws.Range("K99").Value = 1279

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("C6").Value = 10
$ws.Range("K7").Value = 21
$ws.Range("C8").Value = 61
$ws.Range("K8").Value = 40

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("B9").Value = 43
$ws.Range("B10").Value = 64

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("B6").Value = 27
$ws.Range("B7").Value = 42

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("H6").Value = 41
$ws.Range("F7").Value = 71
$ws.Range("F8").Value = 107
$ws.Range("H8").Value = 81

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("G5").Value = 5
$ws.Range("G7").Value = 15

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("F7").Value = 45
$ws.Range("J7").Value = 24
$ws.Range("F9").Value = 159
$ws.Range("J9").Value = 62

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("C9").Value = 36
$ws.Range("C10").Value = 44

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("D3").Value = 7
$ws.Range("D9").Value = 66

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("G3").Value = 15
$ws.Range("J3").Value = 26
$ws.Range("C8").Value = 28
$ws.Range("E8").Value = 53
$ws.Range("B9").Value = 149
$ws.Range("C9").Value = 235
$ws.Range("D9").Value = 410
$ws.Range("E9").Value = 495
$ws.Range("F9").Value = 465
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 86
$ws.Range("B10").Value = 187
$ws.Range("C10").Value = 278
$ws.Range("D10").Value = 472
$ws.Range("E10").Value = 559
$ws.Range("F10").Value = 524
$ws.Range("G10").Value = 213
$ws.Range("J10").Value = 184
$ws.Range("K10").Value = 167

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 3
$ws.Range("F6").Value = 29
$ws.Range("F7").Value = 52
$ws.Range("I7").Value = 18

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("C6").Value = 14
$ws.Range("C7").Value = 30
$ws.Range("C8").Value = 47

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("H6").Value = 5
$ws.Range("H7").Value = 13

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("C7").Value = 7
$ws.Range("C9").Value = 34

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("F5").Value = 34
$ws.Range("F6").Value = 42

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 6
$ws.Range("G3").Value = 6
$ws.Range("C8").Value = 74
$ws.Range("F8").Value = 52
$ws.Range("I8").Value = 37
$ws.Range("C9").Value = 113
$ws.Range("F9").Value = 93
$ws.Range("G9").Value = 70
$ws.Range("I9").Value = 75
$ws.Range("J9").Value = 44

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("H3").Value = 6
$ws.Range("B7").Value = 42
$ws.Range("E7").Value = 47
$ws.Range("B8").Value = 44
$ws.Range("E8").Value = 58
$ws.Range("H8").Value = 33

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("B6").Value = 21
$ws.Range("B7").Value = 22

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("F5").Value = 13
$ws.Range("F6").Value = 23

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("D2").Value = 2
$ws.Range("D8").Value = 31

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("D7").Value = 18
$ws.Range("D8").Value = 21

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 7

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J6").Value = 1
$ws.Range("J8").Value = 10

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("F5").Value = 6
$ws.Range("F6").Value = 9

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("D7").Value = 12
$ws.Range("D8").Value = 15

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("E6").Value = 68
$ws.Range("E7").Value = 73

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 6

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J4").Value = 6
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 14

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 12
